# --- Test Summary sheet removal -------------------------------------------------
$wb = $excel.ActiveWorkbook
$wb.Worksheets.Item("Test Summary").Delete()

$ws = $wb.Worksheets.Item("Test Cases")
$ws.Activate()

# --- Fix the mislabeled TC_UP_00x ids on rows 90-92 (were all "TC_UP_001") ------
$ws.Range("A90").Value = "TC_UP_002"
$ws.Range("A91").Value = "TC_UP_003"
$ws.Range("A92").Value = "TC_UP_004"

# --- Build the new rows 100-105 --------------------------------------------------
# Seed formatting for each new row by copying the closest-matching existing row,
# then patch the handful of cells whose style differs from that template.
$ws.Range("A99:G99").Copy()
$ws.Range("A100:G105").PasteSpecial(-4122)

$ws.Range("B90").Copy()
$ws.Range("B102").PasteSpecial(-4122)
$ws.Range("B104").PasteSpecial(-4122)

$ws.Range("E92").Copy()
$ws.Range("E101").PasteSpecial(-4122)
$ws.Range("E102:E104").PasteSpecial(-4122)

$ws.Range("D78").Copy()
$ws.Range("D105").PasteSpecial(-4122)
$ws.Range("E78").Copy()
$ws.Range("E105").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Row heights as captured in the authored sheet.
$ws.Rows.Item(100).RowHeight = 46.5
$ws.Rows.Item(101).RowHeight = 43.5
$ws.Rows.Item(102).RowHeight = 45.5
$ws.Rows.Item(103).RowHeight = 47.5
$ws.Rows.Item(104).RowHeight = 44
$ws.Rows.Item(105).RowHeight = 46

# Column B (descriptions) filled first across all six new rows ...
$ws.Range("B100").Value = "Check that user is able to navigate through all the products across different categories."
$ws.Range("B101").Value = "Check that all the links are redirecting to correct product/category pages and none of the links are broken."
$ws.Range("B102").Value = "Check that the company logo is clearly visible."
$ws.Range("B103").Value = "Check that all the text – product, category name, price and product description are clearly visible."
$ws.Range("B104").Value = "Check that all the images are clearly visible."
$ws.Range("B105").Value = "Check that category pages have a relevant product listed specifically for the category."

# ... then column A (test case ids) for all six new rows ...
$ws.Range("A100").Value = "TC_M_004"
$ws.Range("A101").Value = "TC_M_005"
$ws.Range("A102").Value = "TC_M_006"
$ws.Range("A103").Value = "TC_M_007"
$ws.Range("A104").Value = "TC_M_008"
$ws.Range("A105").Value = "TC_M_009"

# ... then the expected/remarks columns.
$ws.Range("D100").Value = "user should navigate successfully"
$ws.Range("D101").Value = "should redirect to correct product"
$ws.Range("E101").Value = "redirected to correct product"
$ws.Range("D102").Value = "should be visible"
$ws.Range("E102").Value = "visible"
$ws.Range("D103").Value = "should be visible"
$ws.Range("E103").Value = "visible"
$ws.Range("D104").Value = "should be visible"
$ws.Range("E104").Value = "visible"

# Remaining columns reuse already-existing shared strings.
$ws.Range("C100").Value = "Not applicable"
$ws.Range("C101").Value = "Not applicable"
$ws.Range("C102").Value = "Not applicable"
$ws.Range("C103").Value = "Not applicable"
$ws.Range("C104").Value = "Not applicable"
$ws.Range("C105").Value = "Not applicable"

$ws.Range("E100").Value = "navigated successfully"

$ws.Range("F100").Value = "Pass"
$ws.Range("F101").Value = "Pass"
$ws.Range("F102").Value = "Pass"
$ws.Range("F103").Value = "Pass"
$ws.Range("F104").Value = "Pass"
$ws.Range("F105").Value = "Pass"

# D105/E105 stay blank (only formatting carried over).

# --- View state: scrolled down to show the new rows, B81 selected --------------
$ws.Range("A101").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("B81").Select()
